$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.589.71'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").Value = '2.141.28'
$ws.Range("E3").Value = '  +1.69%  '

$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").Value = "'351.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.14%  '

$ws.Range("D6").Value = "'1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.25%  '

$ws.Range("D7").Value = "'0.5257"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.75%  '

$ws.Range("D8").Value = "'0.4558"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.23%  '

$ws.Range("D9").Value = "'53.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.48%  '

$ws.Range("D10").Value = "'0.09144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.50%  '

$ws.Range("E11").Value = '  +0.39%  '

$ws.Range("D12").Value = "'25.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.99%  '

$ws.Range("D13").Value = '2.138.55'
$ws.Range("E13").Value = '  +1.42%  '

$ws.Range("D14").Value = "'6.878"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.91%  '

$ws.Range("D15").Value = "'8.154"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.86%  '

$ws.Range("D16").Value = "'102.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.52%  '

$ws.Range("E17").Value = '  +2.43%  '

$ws.Range("D18").Value = "'1.009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.22%  '

$ws.Range("D19").Value = "'0.06715"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.07%  '

$ws.Range("D20").Value = "'19.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.91%  '

$ws.Range("D21").Value = "'1.007"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.29%  '

$ws.Range("D22").Value = "'6.347"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.51%  '

$ws.Range("D23").Value = '30.703.34'
$ws.Range("E23").Value = '  +0.57%  '

$ws.Range("D24").Value = "'12.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.06%  '

$ws.Range("D25").Value = "'2.377"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.21%  '

$ws.Range("D26").Value = '2.368.59'
$ws.Range("E26").Value = '  +0.50%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = "'2.661"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.05%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'22.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.64%  '

$ws.Range("D29").Value = "'164.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.20%  '

$ws.Range("D30").Value = "'136.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.14%  '

$ws.Range("E31").Value = '  +1.22%  '

$ws.Range("D33").Value = "'1.669"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.33%  '

$ws.Range("D34").Value = "'6.375"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.22%  '

$ws.Range("D35").Value = "'4.011"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.67%  '

$ws.Range("D36").Value = "'6.177"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.60%  '

$ws.Range("D37").Value = "'10.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.59%  '

$ws.Range("D38").Value = "'0.02653"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.91%  '

$ws.Range("D39").Value = "'0.06948"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.68%  '

$ws.Range("D40").Value = "'0.2339"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.48%  '

$ws.Range("D41").Value = "'12.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.39%  '

$ws.Range("D42").Value = "'0.7021"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.19%  '

$ws.Range("D43").Value = "'1.273"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.11%  '

$ws.Range("D44").Value = "'14.76"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'2.361"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.80%  '

$ws.Range("D46").Value = "'0.6467"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.71%  '

$ws.Range("E47").Value = '  +4.97%  '

$ws.Range("D48").Value = "'3.754"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.49%  '

$ws.Range("D49").Value = "'1.257"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.60%  '

$ws.Range("D50").Value = "'83.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.17%  '

$ws.Range("D51").Value = "'0.07303"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.40%  '
